$wb = $excel.ActiveWorkbook

# Sheet "2025": update row 2 values
$ws = $wb.Worksheets.Item("2025")
$ws.Range("B2").Value = 10372.65132737054
$ws.Range("E2").Value = 289260.5393052954
$ws.Range("G2").Value = 80959.25712661834
$ws.Range("I2").Value = 161710.6685703679
$ws.Range("L2").Value = 484922.2142001599
$ws.Range("M2").Value = 105953.7713982
$ws.Range("N2").Value = 70003.73489578845
$ws.Range("O2").Value = 69744.89343456978

# Sheet "2030": update row 2 values
$ws = $wb.Worksheets.Item("2030")
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 31203.23858116339
$ws.Range("E2").Value = 170658.5511254234
$ws.Range("I2").Value = 209080.6134235085
$ws.Range("L2").Value = 63518.11613148725
$ws.Range("M2").Value = 68536.72857011756
$ws.Range("N2").Value = 19285.19160463996
$ws.Range("O2").Value = 27033.1386905727

# Sheet "2035": update row 2 values
$ws = $wb.Worksheets.Item("2035")
$ws.Range("A2").Value = 27543.1755456332
$ws.Range("B2").Value = 22113.21643273498
$ws.Range("E2").Value = 114655.4402706629
$ws.Range("I2").Value = 153866.0861464091
$ws.Range("M2").Value = 44638.22942194272
$ws.Range("N2").Value = 39676.88529639924
$ws.Range("O2").Value = 31311.04369977792

# Sheet "2040": update row 2 values
$ws = $wb.Worksheets.Item("2040")
$ws.Range("N2").Value = 1142.580190039942
$ws.Range("O2").Value = 0

# Sheet "2045": update row 2 values
$ws = $wb.Worksheets.Item("2045")
$ws.Range("A2").Value = 29588.33508286276
$ws.Range("N2").Value = 4347.543515635315
$ws.Range("O2").Value = 20429.76977394434
